# Add a new "backup" column (R) and append 6 new monthly rows (43-48) to
# the IRFC.NS stock history sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cell R1 = "backup", formatted like the other headers.
# ---------------------------------------------------------------------
$ws.Cells.Item(1, 18).Value = "backup"
$ws.Cells.Item(1, 17).Copy()
$ws.Cells.Item(1, 18).PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 2. Fill the "backup" values for the existing rows (2-42).
# ---------------------------------------------------------------------
$backupValues = @(0,0,0,0,0,0,0,2,1,0,0,0,2,0,0,0,2,0,0,0,0,1,0,0,0,2,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0)
for ($i = 0; $i -lt $backupValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 18).Value = $backupValues[$i]
}

# ---------------------------------------------------------------------
# 3. Append 6 new monthly rows (43-48), matching the formatting of the
#    last existing data row (42). Include column R so the (empty)
#    "backup" cell also exists on the new rows, same as columns A-Q.
# ---------------------------------------------------------------------
$ws.Range("A42:R42").Copy()
$ws.Range("A43:R48").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

$newRows = @(
    @(45474, 172.176549871858, 226.873994153043, 162.6260469672185, 191.8521728515625, 2445064551, 2024, 7, 1, 0, 0, 0, 27, 1, 0, 0),
    @(45505, 192.8428945369496, 193.8336106971541, 173.5041259740911, 176.9617156982422,  644491862, 2024, 8, 1, 0, 0, 0, 31, 0, 0, 0),
    @(45536, 177.8793178319048, 179.2717134191397, 150.8766622945334, 157.8585662841797,  399705004, 2024, 9, 1, 0, 0, 0, 35, 0, 0, 0),
    @(45566, 157.8386597075228, 159.0321437252032, 132.0792368866933, 155.1135406494141,  546024728, 2024, 10, 1, 0, 0, 0, 40, 0, 0, 0),
    @(45597, 156.2970741049755, 160.1261845894605, 137.0521038825714, 148.5294647216797,  374726888, 2024, 11, 1, 0, 0, 0, 44, 0, 0, 2),
    @(45627, 148.8099975585938, 166.8999938964844, 144.6999969482422, 153.3500061035156,  493942055, 2024, 12, 1, 0, 0, 0, 48, 0, 0, 0)
)

# Column order for each data row: A Datetime, B Open, C High, D Low,
# E Close, (F Adj Close left blank), G Volume, H Year, I Month, J Day,
# K Hour, L Minute, M Second, N Week, O isPivot, P two_line_structure,
# Q detect_structure.
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 43 + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    # Column F (Adj Close) intentionally left empty for the new rows.
    $ws.Cells.Item($r, 7).Value = $data[5]
    $ws.Cells.Item($r, 8).Value = $data[6]
    $ws.Cells.Item($r, 9).Value = $data[7]
    $ws.Cells.Item($r, 10).Value = $data[8]
    $ws.Cells.Item($r, 11).Value = $data[9]
    $ws.Cells.Item($r, 12).Value = $data[10]
    $ws.Cells.Item($r, 13).Value = $data[11]
    $ws.Cells.Item($r, 14).Value = $data[12]
    $ws.Cells.Item($r, 15).Value = $data[13]
    $ws.Cells.Item($r, 16).Value = $data[14]
    $ws.Cells.Item($r, 17).Value = $data[15]
    # Column R (backup) intentionally left empty for the new rows.
}
